# Fruta / hortaliza, semanal
# Insert this week's two new price rows (Primera / Segunda) at the top of the
# "Brócoli" data block (row 596), pushing the previously-existing rows down by
# two rows (596-609 -> 598-611), and growing the sheet's used range from
# A1:R609 to A1:R611.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 596 (this shifts the old rows 596-609 down to 598-611)
$ws.Rows.Item(596).Insert()
$ws.Rows.Item(597).Insert()

# New row 596: Brócoli, Primera, Región Metropolitana
$ws.Cells.Item(596,1).Value  = 11
$ws.Cells.Item(596,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(596,3).Value  = "Bíobío"
$ws.Cells.Item(596,4).Value  = 45239
$ws.Cells.Item(596,5).Value  = 8
$ws.Cells.Item(596,6).Value  = 100112023
$ws.Cells.Item(596,7).Value  = "Brócoli"
$ws.Cells.Item(596,8).Value  = "Sin especificar"
$ws.Cells.Item(596,9).Value  = "Primera"
$ws.Cells.Item(596,10).Value = 1000
$ws.Cells.Item(596,11).Value = 1000
$ws.Cells.Item(596,12).Value = 1000
$ws.Cells.Item(596,13).Value = 1000
$ws.Cells.Item(596,14).Value = "$/unidad"
$ws.Cells.Item(596,15).Value = "Región Metropolitana"
$ws.Cells.Item(596,16).Value = 1000
$ws.Cells.Item(596,17).Value = 1
$ws.Cells.Item(596,18).Value = "Hortaliza"

# New row 597: Brócoli, Segunda, Región Metropolitana
$ws.Cells.Item(597,1).Value  = 11
$ws.Cells.Item(597,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(597,3).Value  = "Bíobío"
$ws.Cells.Item(597,4).Value  = 45239
$ws.Cells.Item(597,5).Value  = 8
$ws.Cells.Item(597,6).Value  = 100112023
$ws.Cells.Item(597,7).Value  = "Brócoli"
$ws.Cells.Item(597,8).Value  = "Sin especificar"
$ws.Cells.Item(597,9).Value  = "Segunda"
$ws.Cells.Item(597,10).Value = 1000
$ws.Cells.Item(597,11).Value = 700
$ws.Cells.Item(597,12).Value = 700
$ws.Cells.Item(597,13).Value = 700
$ws.Cells.Item(597,14).Value = "$/unidad"
$ws.Cells.Item(597,15).Value = "Región Metropolitana"
$ws.Cells.Item(597,16).Value = 700
$ws.Cells.Item(597,17).Value = 1
$ws.Cells.Item(597,18).Value = "Hortaliza"
